$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "0.9995", "1.004") are written back as text, matching
# the original inlineStr cell type instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '24.265.10'
$ws.Range('E2').Value = '  +14.39%  '
$ws.Range('D3').Value = '1.672.59'
$ws.Range('E3').Value = '  +8.21%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').Value = '308.23'
$ws.Range('E5').Value = '  +9.13%  '
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('D7').Value = '0.3721'
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('D8').Value = '0.3422'
$ws.Range('E8').Value = '  +6.86%  '
$ws.Range('D9').Value = '47.64'
$ws.Range('E9').Value = '  +16.21%  '
$ws.Range('D10').Value = '1.176'
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('D11').Value = '0.07260'
$ws.Range('E11').Value = '  +6.02%  '
$ws.Range('D12').Value = '0.9998'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '20.40'
$ws.Range('E13').Value = '  +8.42%  '
$ws.Range('D14').Value = '6.072'
$ws.Range('E14').Value = '  +6.34%  '
$ws.Range('D15').Value = '6.735'
$ws.Range('E15').Value = '  +5.45%  '
$ws.Range('D16').Value = '1.674.80'
$ws.Range('E16').Value = '  +8.88%  '
$ws.Range('D17').Value = '0.00001104'
$ws.Range('E17').Value = '  +5.43%  '
$ws.Range('D18').Value = '0.9991'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').Value = '0.06702'
$ws.Range('E19').Value = '  +9.16%  '
$ws.Range('D20').Value = '81.26'
$ws.Range('E20').Value = '  +11.54%  '
$ws.Range('D21').Value = '16.39'
$ws.Range('E21').Value = '  +8.25%  '
$ws.Range('D22').Value = '6.119'
$ws.Range('E22').Value = '  +6.65%  '
$ws.Range('D23').Value = '12.02'
$ws.Range('E23').Value = '  +5.61%  '
$ws.Range('D24').Value = '24.199.47'
$ws.Range('E24').Value = '  +13.82%  '
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('D26').Value = '3.359'
$ws.Range('E26').Value = '  -9.27%  '
$ws.Range('D27').Value = '2.637'
$ws.Range('E27').Value = '  +17.64%  '
$ws.Range('D28').Value = '152.02'
$ws.Range('E28').Value = '  +2.18%  '
$ws.Range('D29').Value = '19.40'
$ws.Range('E29').Value = '  +9.39%  '
$ws.Range('D30').Value = '1.861.40'
$ws.Range('E30').Value = '  +8.96%  '
$ws.Range('D31').Value = '126.88'
$ws.Range('E31').Value = '  +6.85%  '
$ws.Range('D32').Value = '6.386'
$ws.Range('E32').Value = '  +21.69%  '
$ws.Range('D33').Value = '4.044'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('D34').Value = '0.9848'
$ws.Range('E34').Value = '  +14.25%  '
$ws.Range('D35').Value = '1.748'
$ws.Range('E35').Value = '  +15.82%  '
$ws.Range('D36').Value = '0.08441'
$ws.Range('E36').Value = '  +5.06%  '
$ws.Range('D37').Value = '12.52'
$ws.Range('E37').Value = '  +16.75%  '
$ws.Range('D38').Value = '0.06407'
$ws.Range('E38').Value = '  +8.97%  '
$ws.Range('D39').Value = '5.349'
$ws.Range('E39').Value = '  +7.84%  '
$ws.Range('D40').Value = '8.763'
$ws.Range('E40').Value = '  +11.82%  '
$ws.Range('D41').Value = '0.02334'
$ws.Range('E41').Value = '  +10.42%  '
$ws.Range('D42').Value = '1.280'
$ws.Range('E42').Value = '  +5.49%  '
$ws.Range('D43').Value = '0.2101'
$ws.Range('E43').Value = '  +9.10%  '
$ws.Range('D44').Value = '0.6138'
$ws.Range('E44').Value = '  +12.08%  '
$ws.Range('D45').Value = '0.9987'
$ws.Range('E45').Value = '  +3.32%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.802'
$ws.Range('E46').Value = '  +6.22%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '13.10'
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('D48').Value = '0.5931'
$ws.Range('E48').Value = '  +8.59%  '
$ws.Range('D49').Value = '127.77'
$ws.Range('E49').Value = '  +4.59%  '
$ws.Range('D50').Value = '2.011'
$ws.Range('E50').Value = '  +6.95%  '
$ws.Range('E51').Value = '  +7.76%  '
